$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.583.26'
$ws.Range('E2').Value = '  -1.62%  '
$ws.Range('D3').Value = '2.064.85'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  -2.49%  '
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''52.69'
$ws.Range('E8').Value = '  -8.04%  '
$ws.Range('D9').Value = '''59.11'
$ws.Range('E9').Value = '  -1.86%  '
$ws.Range('D10').Value = '''0.359'
$ws.Range('E10').Value = '  -7.09%  '
$ws.Range('D11').Value = '''0.0750'
$ws.Range('E11').Value = '  -4.59%  '
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').Value = '''0.906'
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('D14').Value = '''14.65'
$ws.Range('E14').Value = '  -9.78%  '
$ws.Range('D15').Value = '2.363.75'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '''5.40'
$ws.Range('E16').Value = '  -6.05%  '
$ws.Range('D17').Value = '2.086.70'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '36.496.84'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('D19').Value = '''16.37'
$ws.Range('E19').Value = '  -12.87%  '
$ws.Range('D20').Value = '''71.73'
$ws.Range('E20').Value = '  -4.28%  '
$ws.Range('E21').Value = '  -4.19%  '
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '''5.25'
$ws.Range('E23').Value = '  -4.19%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E25').Value = '  -4.84%  '
$ws.Range('D26').Value = '''9.43'
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('D28').Value = '''164.21'
$ws.Range('E28').Value = '  -3.76%  '
$ws.Range('D29').Value = '''20.53'
$ws.Range('E29').Value = '  +1.48%  '
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').Value = '''5.04'
$ws.Range('E31').Value = '  -2.12%  '
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('D33').Value = '''4.56'
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('D34').Value = '''0.0596'
$ws.Range('E34').Value = '  -4.49%  '
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''1.84'
$ws.Range('E36').Value = '  +3.58%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '''2.27'
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('D38').Value = '''0.0817'
$ws.Range('E38').Value = '  -7.39%  '
$ws.Range('D39').Value = '''1.25'
$ws.Range('E39').Value = '  -7.00%  '
$ws.Range('D40').Value = '''2.93'
$ws.Range('E40').Value = '  -4.58%  '
$ws.Range('D41').Value = '''4.85'
$ws.Range('E41').Value = '  -7.36%  '
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('E43').Value = '  -3.82%  '
$ws.Range('D44').Value = '''0.0935'
$ws.Range('E44').Value = '  -6.96%  '
$ws.Range('D45').Value = '''94.21'
$ws.Range('D46').Value = '1.390.65'
$ws.Range('E46').Value = '  +8.74%  '
$ws.Range('D47').Value = '''7.40'
$ws.Range('E47').Value = '  +7.95%  '
$ws.Range('D48').Value = '''15.50'
$ws.Range('E48').Value = '  -11.69%  '
$ws.Range('D49').Value = '''2.34'
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').Value = '2.252.83'
$ws.Range('E51').Value = '  +0.20%  '
